$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: new weekly date headers, newest-first, old ones shift right ---
# Write the nine brand-new dates first, oldest-to-newest (matches shared-string
# allocation order), then the carried-forward dates (already in the shared
# string table, order doesn't matter for allocation).
$newDates = @("Jun_16","Jun_24","Jun_30","Jul_07","Jul_17","Jul_23","Aug_04","Aug_25","Sep_08")
for ($i = 0; $i -lt $newDates.Length; $i++) {
    # newest date -> column B (2); oldest of the new dates -> column J (10)
    $col = 10 - $i
    $ws.Cells.Item(1, $col).Value = $newDates[$i]
}

$carried = @("Jun_09","Jun_03","May_27","May_23","May_19","May_15","May_12","May_05","Apr_28","Apr_24","Apr_21","Apr_17","Apr_11")
for ($i = 0; $i -lt $carried.Length; $i++) {
    $ws.Cells.Item(1, 11 + $i).Value = $carried[$i]
}

# --- Rows 2-31: extend existing "UN" run from column B through column W ---
for ($r = 2; $r -le 31; $r++) {
    for ($c = 15; $c -le 23; $c++) {
        $ws.Cells.Item($r, $c).Value = "UN"
    }
}

# --- Rows 32-33: extend existing "UN" run from column B through column P ---
for ($r = 32; $r -le 33; $r++) {
    for ($c = 8; $c -le 16; $c++) {
        $ws.Cells.Item($r, $c).Value = "UN"
    }
}
